# "Cambié el orden de unas diapos" — reorder two slides (positions 15 and 17)
# and carry out the small content tweaks that came along with that move.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Reorder: swap the slides that sit at position 15 and position 17
#    (the slide that used to be in the middle, position 16, stays put).
# ---------------------------------------------------------------------------
$p.Slides.Item(17).MoveTo(15)
$p.Slides.Item(16).MoveTo(17)

# ---------------------------------------------------------------------------
# 2) The slide now sitting at position 15 (previously at 17): it had an
#    empty "Title 1" placeholder, an empty "Content Placeholder 2" and a
#    picture. Drop the two empty placeholders and replace the title with a
#    normal textbox that reads "Uso en Chile".
# ---------------------------------------------------------------------------
$s15 = $p.Slides.Item(15)

# Remove the (empty) content placeholder.
$s15.Shapes.Item("Content Placeholder 2").Delete()

# Remove the (empty) title placeholder. The first Delete() on this
# particular placeholder only clears it back to its empty layout default,
# so it has to be called twice to actually get rid of the shape.
$s15.Shapes.Item("Title 1").Delete()
$s15.Shapes.Item("Title 1").Delete()

# Add the freestanding title textbox that replaces it.
$tb = $s15.Shapes.AddTextbox(1, 32.24007874015748, 7.713779527559055, 701.7076377952756, 100.85748031496063)
$tb.Name = "Title 1"
$tb.TextFrame.AutoSize = 2
$tb.TextFrame.MarginLeft = 7.2
$tb.TextFrame.MarginRight = 7.2
$tb.TextFrame.MarginTop = 3.6
$tb.TextFrame.MarginBottom = 3.6
$tb.TextFrame.VerticalAnchor = 1
$tb.TextFrame.TextRange.Text = "Uso en Chile"

# ---------------------------------------------------------------------------
# 3) The slide now sitting at position 17 (previously at 15): it had the
#    "Uso en Chile" title and a picture placeholder. Drop the title and
#    enlarge/reposition the picture to fill the freed-up space.
# ---------------------------------------------------------------------------
$s17 = $p.Slides.Item(17)

$s17.Shapes.Item("Title 1").Delete()
$s17.Shapes.Item("Title 1").Delete()

$pic17 = $s17.Shapes.Item("Picture 2")
$pic17.Left = 154.37181102362206
$pic17.Top = 54
$pic17.Width = 651.2563779527559
$pic17.Height = 432

# ---------------------------------------------------------------------------
# 4) Slide 8: nudge the first picture a bit to the right/up.
# ---------------------------------------------------------------------------
$s8 = $p.Slides.Item(8)
$pic8 = $s8.Shapes.Item("Picture 2")
$pic8.Left = 110.91629921259843
$pic8.Top = 157.93425196850393
